{"js": "const body = context.document.body;\n\n// --- Edit 1 -------------------------------------------------------------\n// Original run: \" We then summed up the migration from Democratic states and\n// Republican states. \" (trailing \". \" at the end).\n// New text:     \" We then summed up the migration from Democratic states and\n// Republican states into a variable called totalMigration. \"\n// i.e. the trailing \". \" is replaced with \" into a variable called\n// totalMigration. \".\nconst states = body.search(\"Republican states. \", { matchCase: true });\nstates.load(\"items\");\nawait context.sync();\n\nif (states.items.length !== 1) {\n  throw new Error(\"Expected exactly 1 match for 'Republican states. ', got \" + states.items.length);\n}\nstates.items[0].insertText(\n  \"Republican states into a variable called totalMigration. \",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Edit 2 -------------------------------------------------------------\n// Right after \"...a negative value would mean Republican migration. \" add a\n// new sentence: \" We did this for 2016 and 2020, so we had a total skew for\n// both elections.\"\nconst skew = body.search(\n  \"a negative value would mean Republican migration. \",\n  { matchCase: true }\n);\nskew.load(\"items\");\nawait context.sync();\n\nif (skew.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for 'a negative value would mean Republican migration. ', got \" +\n      skew.items.length\n  );\n}\nskew.items[0].insertText(\n  \" We did this for 2016 and 2020, so we had a total skew for both elections.\",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------------\n# Original run ends \"...Democratic states and Republican states. \" (with a\n# trailing \". \"). Replace that trailing \". \" with\n# \" into a variable called totalMigration. \" so the sentence reads\n# \"...Republican states into a variable called totalMigration. \".\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\n    \"Republican states. \",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"Republican states into a variable called totalMigration. \",\n    2\n)\nif (-not $found1) {\n    throw \"Could not find 'Republican states. ' to replace.\"\n}\n\n# --- Edit 2 -----------------------------------------------------------------\n# Append a new sentence right after \"...Republican migration. \":\n# \" We did this for 2016 and 2020, so we had a total skew for both elections.\"\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\n    \"a negative value would mean Republican migration. \",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"a negative value would mean Republican migration.  We did this for 2016 and 2020, so we had a total skew for both elections.\",\n    2\n)\nif (-not $found2) {\n    throw \"Could not find 'a negative value would mean Republican migration. ' to replace.\"\n}\n"}
